# Apply the "Add repost, fix excel google sheets" change:
# - Insert a new first column (shifts id.. Блокировка бота from A:J to B:K)
# - Set new header "№" in A1
# - Add two data rows (2 and 3) under the header
# - Leave row 4 empty
# - Add a second, smaller header block on row 5 and a data row on row 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the left, shifting existing columns A:J to B:K
$ws.Columns.Item(1).Insert()

# New header cell for the inserted column
$ws.Cells.Item(1, 1).Value = "№"

# Row 2 - first data row
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 1012882762
$ws.Cells.Item(2, 3).Value = "alekseinushtaev"
$ws.Cells.Item(2, 4).Value = "Алексей"
$ws.Cells.Item(2, 6).Value = "2025-04-11   14:25:18"
$ws.Cells.Item(2, 7).Value = "test"
$ws.Cells.Item(2, 8).Value = "test"
$ws.Cells.Item(2, 9).Value = "test"
$ws.Cells.Item(2, 10).Value = "2025-04-12   10:54:46"

# Row 3 - second data row
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 1012882762
$ws.Cells.Item(3, 3).Value = "alekseinushtaev"
$ws.Cells.Item(3, 4).Value = "Алексей"
$ws.Cells.Item(3, 6).Value = "2025-04-11   14:25:18"
$ws.Cells.Item(3, 7).Value = "test2"
$ws.Cells.Item(3, 8).Value = "test2"
$ws.Cells.Item(3, 9).Value = "test2"
$ws.Cells.Item(3, 10).Value = "2025-04-12   10:54:58"

# Row 5 - second (shorter) header block
$ws.Cells.Item(5, 1).Value = "№"
$ws.Cells.Item(5, 2).Value = "id"
$ws.Cells.Item(5, 3).Value = "username"
$ws.Cells.Item(5, 4).Value = "first_name"
$ws.Cells.Item(5, 5).Value = "last_name"
$ws.Cells.Item(5, 6).Value = "Время входа в бота"
$ws.Cells.Item(5, 7).Value = "Блокировка бота"

# Row 6 - data row for the second block
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = 1012882762
$ws.Cells.Item(6, 3).Value = "alekseinushtaev"
$ws.Cells.Item(6, 4).Value = "Алексей"
$ws.Cells.Item(6, 6).Value = "2025-04-11   14:25:18"

Write-Host "applied edits"
